# Insert two new price-record rows for "Naranja" (Fukumoto / Valencia) at
# the top of the weekly block that starts at row 700, pushing all the
# existing rows (700..755) down to (702..757).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 700.
$ws.Rows("700:701").Insert()

# --- New row 700: Naranja / Fukumoto / Primera -----------------------------
$ws.Range("A700").Value = 9
$ws.Range("B700").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C700").Value = "Metropolitana"
$ws.Range("D700").Value = 44714
$ws.Range("E700").Value = 13
$ws.Range("F700").Value = "Fruta"
$ws.Range("G700").Value = 100102
$ws.Range("H700").Value = "Cítricos"
$ws.Range("I700").Value = 100102005
$ws.Range("J700").Value = "Naranja"
$ws.Range("K700").Value = "Fukumoto"
$ws.Range("L700").Value = "Primera"
$ws.Range("M700").Value = 300
$ws.Range("N700").Value = 11000
$ws.Range("O700").Value = 11000
$ws.Range("P700").Value = 11000
$ws.Range("Q700").Value = "`$/caja 18 kilos granel"
$ws.Range("R700").Value = "Región de O'Higgins"
$ws.Range("S700").Value = 611
$ws.Range("T700").Value = 18

# --- New row 701: Naranja / Valencia / Primera ------------------------------
$ws.Range("A701").Value = 9
$ws.Range("B701").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C701").Value = "Metropolitana"
$ws.Range("D701").Value = 44714
$ws.Range("E701").Value = 13
$ws.Range("F701").Value = "Fruta"
$ws.Range("G701").Value = 100102
$ws.Range("H701").Value = "Cítricos"
$ws.Range("I701").Value = 100102005
$ws.Range("J701").Value = "Naranja"
$ws.Range("K701").Value = "Valencia"
$ws.Range("L701").Value = "Primera"
$ws.Range("M701").Value = 350
$ws.Range("N701").Value = 10500
$ws.Range("O701").Value = 10500
$ws.Range("P701").Value = 10500
$ws.Range("Q701").Value = "`$/caja 15 kilos granel"
$ws.Range("R701").Value = "Provincia de Melipilla"
$ws.Range("S701").Value = 700
$ws.Range("T701").Value = 15
